$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds a 2-row table (rows 2-3, below the header in row 1)
# describing Narayan Jagadeesan's innings. The update appends two more rows
# (4-5) that repeat the existing row-2 and row-3 innings records respectively.
# Copy/PasteSpecial (rather than re-typing the values) preserves the original
# cell typing (text, even for numeric-looking values like "0", "33", "117.85")
# and exact characters (e.g. the non-breaking space in "Narayan Jagadeesan ").

$ws.Range("A2:K3").Copy()
$ws.Range("A4:K5").PasteSpecial()
